# Add "Romania" and "Slovakia" test-data worksheets, each a copy of the
# existing "Spain" sheet with the market name (B2) and ticket reference
# (B4) updated for the new market.

$wb = $excel.ActiveWorkbook
$spain = $wb.Worksheets.Item("Spain")

# --- Romania sheet: copy of Spain, placed after the last sheet ---
$spain.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$romania = $wb.Worksheets.Item($wb.Worksheets.Count)
$romania.Name = "Romania"
$romania.Range("B2").Value = "Romania Market"
$romania.Range("B4").Value = "NGC-4307/T3537/T3551"

# Spain is no longer the active tab; restore its sheet-level selection to
# the "select all" state left behind once another sheet becomes active.
$spain.Cells.Select() | Out-Null

# --- Slovakia sheet: copy of Romania, placed after the last sheet ---
$romania.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"
$slovakia.Range("B4").Value = "NGC-4306/T3564/T3576"
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B3").Select() | Out-Null
